$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep plain-text storage (matching original inlineStr/text cells)
# by explicitly setting the Text number format before assigning the new value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "321.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.97%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "42.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-5.85%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.148"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-8.40%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08185"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.83%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.282"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.52%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.798"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-14.22%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9325"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.10%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1109"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.50%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1866"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.08%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09443"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.90%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04644"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.70%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.400"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-28.57%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.21%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001302"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.83%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005755"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.59%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.365"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.34%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.528"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.06%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.01%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.27%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2522"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-12.43%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04169"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.24%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001247"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.37%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004363"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.26%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-7.88%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002981"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.51%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02744"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "1.16%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05575"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.35%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.008062"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.59%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1397"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.65%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006547"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-9.94%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002094"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.35%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007561"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-17.09%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3497"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.35%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006958"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.18%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.22%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003478"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.37%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003533"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.71%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.22%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.22%"
